$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.320.96'
$ws.Range('E2').Value = '  +6.06%  '
$ws.Range('D3').Value = '3.350.61'
$ws.Range('E3').Value = '  +2.77%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '412.31'
$ws.Range('E5').Value = '  +3.89%  '
$ws.Range('D6').Value = '112.65'
$ws.Range('E6').Value = '  +2.76%  '
$ws.Range('E7').Value = '  +4.97%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.636'
$ws.Range('E9').Value = '  +2.10%  '
$ws.Range('D10').Value = '40.21'
$ws.Range('E10').Value = '  +2.53%  '
$ws.Range('E11').Value = '  +2.02%  '
$ws.Range('E12').Value = '  +1.39%  '
$ws.Range('D13').Value = '3.879.20'
$ws.Range('E13').Value = '  +2.90%  '
$ws.Range('D14').Value = '8.57'
$ws.Range('D15').Value = '19.40'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').Value = '3.346.03'
$ws.Range('E16').Value = '  +2.71%  '
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range('D18').Value = '60.153.58'
$ws.Range('E18').Value = '  +6.01%  '
$ws.Range('D19').Value = '10.82'
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').Value = '3.37'
$ws.Range('E20').Value = '  +2.15%  '
$ws.Range('E21').Value = '  +4.37%  '
$ws.Range('D22').Value = '13.14'
$ws.Range('E22').Value = '  +2.13%  '
$ws.Range('D23').Value = '306.00'
$ws.Range('D24').Value = '75.89'
$ws.Range('E24').Value = '  +1.27%  '
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '28.70'
$ws.Range('E26').Value = '  +1.84%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '0.183'
$ws.Range('E27').Value = '  +8.78%  '
$ws.Range('D28').Value = '4.48'
$ws.Range('E28').Value = '  +2.37%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.79'
$ws.Range('E29').Value = '  +30.36%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '7.92'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').Value = '7.48'
$ws.Range('E31').Value = '  +3.07%  '
$ws.Range('E32').Value = '  +4.74%  '
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('D34').Value = '11.60'
$ws.Range('E34').Value = '  +5.49%  '
$ws.Range('D35').Value = '39.96'
$ws.Range('E35').Value = '  +6.99%  '
$ws.Range('E36').Value = '  +5.59%  '
$ws.Range('D37').Value = '52.15'
$ws.Range('E37').Value = '  +1.16%  '
$ws.Range('D38').Value = '3.18'
$ws.Range('E38').Value = '  +1.13%  '
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').Value = '3.41'
$ws.Range('E40').Value = '  -4.19%  '
$ws.Range('D41').Value = '138.06'
$ws.Range('E41').Value = '  +2.71%  '
$ws.Range('E42').Value = '  +2.87%  '
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').Value = '3.95'
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '0.285'
$ws.Range('E45').Value = '  +1.64%  '
$ws.Range('D46').Value = '16.96'
$ws.Range('E46').Value = '  -2.25%  '
$ws.Range('E47').Value = '  +8.77%  '
$ws.Range('D48').Value = '22.29'
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('D49').Value = '2.208.43'
$ws.Range('E49').Value = '  +2.63%  '
$ws.Range('D50').Value = '2.06'
$ws.Range('E50').Value = '  +0.59%  '
$ws.Range('D51').Value = '2.41'
$ws.Range('E51').Value = '  +1.26%  '
